# Revert "Drop in all data files from 3.0 RMI script"
# - Insert a new "Texas Data" worksheet between "IEA Data" and "HPEbP"
#   containing a set of explanatory notes about a calculation error found
#   in the NREL-derived efficiency numbers.
# - Fix the natural gas reforming efficiency formula on the HPEbP sheet
#   (cell B3) to stop counting waste heat as an energy input.

$wb = $excel.ActiveWorkbook

# --- 1. Insert the "Texas Data" worksheet right after "IEA Data" ---------
$iea = $wb.Worksheets.Item("IEA Data")
$tx = $wb.Worksheets.Add($null, $iea)
$tx.Name = "Texas Data"

$tx.Range("A1").Value = "There is no reason that these number should be different for Texas."
$tx.Range("A3").Value = "However, I did find an error in their calculations. "
$tx.Range("A5").Value = "They were included waste heat as an energy balance input."
$tx.Range("B6").Value = "for example, page 228 of the NREL report shows gas production as 162 kBtu gas + 2 kBtu electricity = 118 kBtu hydrogen + 46 kBtu waste heat"
$tx.Range("B7").Value = "so, the efficiency (output hydrogen energy vs input energy) would be 118/(162+2)=72%"
$tx.Range("B8").Value = "previously, this spreadhseet (cell 'HPEbP'B3) was calculating the efficiency as 118/(162+2+46)=56%"
$tx.Range("B10").Value = "the IEA number for natural gas reforming efficiency is 76%, so that's a good check that their initial calculation was wrong. "
$tx.Range("A12").Value = "Their other calculations did not include the same mistake."

# Style every cell in A1:I17 with the "accent5" themed font colour used
# throughout these notes.
$tx.Range("A1:I17").Font.ThemeColor = 9

# --- 2. Fix the HPEbP natural-gas reforming efficiency formula -----------
$hp = $wb.Worksheets.Item("HPEbP")
$hp.Range("B3").Formula = "=118/(162+2)"

# --- 3. Update saved cursor/selection position on each sheet -------------
$about = $wb.Worksheets.Item("About")
$about.Range("B14").Select() | Out-Null

$iea.Range("E18").Select() | Out-Null

$tx.Range("A13").Select() | Out-Null

$hp.Activate()
$hp.Range("C12").Select() | Out-Null
